$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column index order corresponding to each data row's values after the row number:
# E(5), G(7), H(8), I(9), J(10), K(11), M(13), N(14), O(15), P(16), Q(17), R(18), S(19), T(20)
$colIndexes = @(5, 7, 8, 9, 10, 11, 13, 14, 15, 16, 17, 18, 19, 20)

$data = @(
    @(2, 3, 0.7922663333333334, 2.376799, 0.1759587713796512, 0.1759587713796512, 3, 5.760626333333334, 17.281879, 0.1124958508646345, 0.1124958508646345, 4.563950302813445, 41.075552725321, 0.01979463170344956, 0.01979463170344956),
    @(3, 3, 0.7922663333333334, 2.376799, 0.1759587713796512, 0.1759587713796512, 3, 28.56714433333333, 85.70143299999999, 0.5578707978254833, 0.5578707978254833, 22.63278669477411, 203.695080252967, 0.09816226017395782, 0.09816226017395785),
    @(4, 3, 0.7922663333333334, 2.376799, 0.1759587713796512, 0.1759587713796512, 3, 12.81995433333333, 38.459863, 0.2503532753771899, 0.2503532753771899, 10.15681821317078, 91.41136391853701, 0.04405185474624181, 0.04405185474624182),
    @(5, 3, 0.7922663333333334, 2.376799, 0.1759587713796512, 0.1759587713796512, 3, 4.059731, 12.179193, 0.0792800759326923, 0.07928007593269232, 3.216388193689667, 28.947493743207, 0.01395002475600199, 0.013950024756002),
    @(6, 3, 2.077831666666667, 6.233495, 0.4614770208171574, 0.4614770208171574, 3, 5.760626333333334, 17.281879, 0.1124958508646345, 0.1124958508646345, 11.96961181523389, 107.726506337105, 0.05191425011130277, 0.05191425011130277),
    @(7, 3, 2.077831666666667, 6.233495, 0.4614770208171574, 0.4614770208171574, 3, 28.56714433333333, 85.70143299999999, 0.5578707978254833, 0.5578707978254833, 59.35771712203723, 534.2194540983349, 0.2574445537813947, 0.2574445537813948),
    @(8, 3, 2.077831666666667, 6.233495, 0.4614770208171574, 0.4614770208171574, 3, 12.81995433333333, 38.459863, 0.2503532753771899, 0.2503532753771899, 26.63770707902056, 239.739363711185, 0.115532283672883, 0.115532283672883),
    @(9, 3, 2.077831666666667, 6.233495, 0.4614770208171574, 0.4614770208171574, 3, 4.059731, 12.179193, 0.0792800759326923, 0.07928007593269232, 8.435437629948336, 75.91893866953501, 0.03658593325157686, 0.03658593325157687),
    @(10, 3, 0.8480786666666668, 2.544236, 0.1883544383264543, 0.1883544383264543, 3, 5.760626333333334, 17.281879, 0.1124958508646345, 0.1124958508646345, 4.885464299938223, 43.969178699444, 0.0211890928036648, 0.0211890928036648),
    @(11, 3, 0.8480786666666668, 2.544236, 0.1883544383264543, 0.1883544383264543, 3, 28.56714433333333, 85.70143299999999, 0.5578707978254833, 0.5578707978254833, 24.22718567668756, 218.044671090188, 0.1050774407831498, 0.1050774407831498),
    @(12, 3, 0.8480786666666668, 2.544236, 0.1883544383264543, 0.1883544383264543, 3, 12.81995433333333, 38.459863, 0.2503532753771899, 0.2503532753771899, 10.87232977774089, 97.850967999668, 0.04715515056685873, 0.04715515056685874),
    @(13, 3, 0.8480786666666668, 2.544236, 0.1883544383264543, 0.1883544383264543, 3, 4.059731, 12.179193, 0.0792800759326923, 0.07928007593269232, 3.442971253505334, 30.98674128154801, 0.0149327541727809, 0.01493275417278091),
    @(14, 3, 0.7843913333333336, 2.353174000000001, 0.1742097694767371, 0.1742097694767372, 3, 5.760626333333334, 17.281879, 0.1124958508646345, 0.1124958508646345, 4.518585370438446, 40.66726833394601, 0.01959787624621738, 0.01959787624621738),
    @(15, 3, 0.7843913333333336, 2.353174000000001, 0.1742097694767371, 0.1742097694767372, 3, 28.56714433333333, 85.70143299999999, 0.5578707978254833, 0.5578707978254833, 22.40782043314912, 201.670383898342, 0.09718654308698088, 0.09718654308698089),
    @(16, 3, 0.7843913333333336, 2.353174000000001, 0.1742097694767371, 0.1742097694767372, 3, 12.81995433333333, 38.459863, 0.2503532753771899, 0.2503532753771899, 10.05586107279578, 90.50274965516202, 0.04361398639120634, 0.04361398639120635),
    @(17, 3, 0.7843913333333336, 2.353174000000001, 0.1742097694767371, 0.1742097694767372, 3, 4.059731, 12.179193, 0.0792800759326923, 0.07928007593269232, 3.184417812064668, 28.65976030858201, 0.01381136375233254, 0.01381136375233255)

)

foreach ($entry in $data) {
    $r = $entry[0]
    for ($i = 0; $i -lt $colIndexes.Count; $i++) {
        $col = $colIndexes[$i]
        $val = $entry[$i + 1]
        $ws.Cells.Item($r, $col).Value = $val
    }
}

Write-Output "Updated $($data.Count) rows across $($colIndexes.Count) columns."

